$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new data row (row 10) by cloning the formatting of the last existing
# data row (row 9), then overwriting its values with the new review record.
# ---------------------------------------------------------------------------
$ws.Range("A9:G9").Copy() | Out-Null
$ws.Range("A10:G10").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = "com.hamxa.shaynachim"
$ws.Range("B10").Value = "bitcoin"
$ws.Range("C10").Value = "shmulmaor2@gmail.com"
$ws.Range("D10").Value = "vikicrestina@gmail.com"
$ws.Range("E10").Value = "27/5/2019 15:59"
$ws.Range("F10").Value = "I read everything...valuable information guide. A lot to learn"
$ws.Range("G10").Value = "no"

# New hyperlinks for the email columns of the new row (this also applies
# Excel's built-in "Hyperlink" font style to the cells).
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com") | Out-Null

# ---------------------------------------------------------------------------
# Restore the plain bordered cell style on C10:D10 (overwritten above by the
# hyperlink auto-style) and extend that same style down through the new,
# otherwise empty, trailing rows 11-25.
# ---------------------------------------------------------------------------
$ws.Range("C9:D9").Copy() | Out-Null
$ws.Range("C10:D25").PasteSpecial(-4122) | Out-Null

for ($r = 11; $r -le 25; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# ---------------------------------------------------------------------------
# Update the view: scroll back to the left edge and select A11.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A11").Select() | Out-Null
